$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}

Replace-Text "2023-09-28 Thursday" "2023-09-29 Friday"

Replace-Text "79×88=6952" "53×18=954"
Replace-Text "13×90=1170" "45×25=1125"
Replace-Text "66×84=5544" "76×18=1368"
Replace-Text "76×70=5320" "72×77=5544"
Replace-Text "95×89=8455" "43×72=3096"

Replace-Text "95×15=1425" "75×84=6300"
Replace-Text "56×20=1120" "38×61=2318"
Replace-Text "67×46=3082" "70×47=3290"
Replace-Text "36×80=2880" "56×61=3416"
Replace-Text "89×18=1602" "78×74=5772"

Replace-Text "40×41=1640" "61×83=5063"
Replace-Text "53×47=2491" "72×63=4536"
Replace-Text "46×88=4048" "72×22=1584"
Replace-Text "91×46=4186" "31×46=1426"
Replace-Text "52×17=884" "54×46=2484"

Replace-Text "80×85=6800" "35×20=700"
Replace-Text "26×83=2158" "38×56=2128"
Replace-Text "39×86=3354" "91×21=1911"
Replace-Text "19×39=741" "72×17=1224"
Replace-Text "81×73=5913" "99×70=6930"

Replace-Text "27×19=513" "49×16=784"
Replace-Text "36×38=1368" "42×42=1764"
Replace-Text "82×49=4018" "28×89=2492"
Replace-Text "49×22=1078" "57×91=5187"
Replace-Text "55×74=4070" "87×18=1566"
